$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update country labels (reordering of countries by updated case totals)
$ws.Range("A14").Value = "Brasil"
$ws.Range("A15").Value = "Canada"
$ws.Range("A16").Value = "Paises Bajos"
$ws.Range("A17").Value = "Suiza"
$ws.Range("A193").Value = "Republica de Africa Central"
$ws.Range("A194").Value = "Seychelles"
$ws.Range("A195").Value = "Islas Malvinas"
$ws.Range("A196").Value = "Montserrat"

# Update statistic values (Casos totales, Nuevos casos, Casos activos, Recuperados,
# Casos criticos, Muertes hoy, Muertes) for the affected countries

# Row 4: Estados Unidos
$ws.Range("B4").Value = 623694
$ws.Range("C4").Value = 9808
$ws.Range("E4").Value = 548382
$ws.Range("G4").Value = 1558
$ws.Range("H4").Value = 27605

# Row 8: Alemania
$ws.Range("B8").Value = 133456
$ws.Range("C8").Value = 1246
$ws.Range("E8").Value = 57264

# Row 14: now Brasil
$ws.Range("B14").Value = 28320
$ws.Range("C14").Value = 3058
$ws.Range("D14").Value = 14026
$ws.Range("E14").Value = 12558
$ws.Range("F14").Value = 296
$ws.Range("G14").Value = 204
$ws.Range("H14").Value = 1736

# Row 15: now Canada
$ws.Range("B15").Value = 28205
$ws.Range("C15").Value = 1142
$ws.Range("D15").Value = 8937
$ws.Range("E15").Value = 18262
$ws.Range("F15").Value = 557
$ws.Range("G15").Value = 103
$ws.Range("H15").Value = 1006

# Row 16: now Paises Bajos
$ws.Range("B16").Value = 28153
$ws.Range("C16").Value = 734
$ws.Range("D16").Value = 250
$ws.Range("E16").Value = 24769
$ws.Range("F16").Value = 1279
$ws.Range("G16").Value = 189
$ws.Range("H16").Value = 3134

# Row 17: now Suiza
$ws.Range("B17").Value = 26336
$ws.Range("C17").Value = 400
$ws.Range("D17").Value = 15400
$ws.Range("E17").Value = 9697
$ws.Range("F17").Value = 386
$ws.Range("G17").Value = 65
$ws.Range("H17").Value = 1239

# Row 20: Austria
$ws.Range("B20").Value = 14336
$ws.Range("C20").Value = 110
$ws.Range("E20").Value = 5845

# Row 93: Costa Rica
$ws.Range("B93").Value = 626
$ws.Range("C93").Value = 8
$ws.Range("D93").Value = 67
$ws.Range("E93").Value = 555
$ws.Range("F93").Value = 11
$ws.Range("G93").Value = 1
$ws.Range("H93").Value = 4

# Row 193: now Republica de Africa Central
$ws.Range("B193").Value = 12
$ws.Range("C193").Value = 1
$ws.Range("D193").Value = 4
$ws.Range("E193").Value = 8

# Row 194: now Seychelles
$ws.Range("D194").Value = 0
$ws.Range("E194").Value = 11

# Row 195: now Islas Malvinas
$ws.Range("F195").Value = 0

# Row 196: now Montserrat
$ws.Range("D196").Value = 1
$ws.Range("E196").Value = 10
$ws.Range("F196").Value = 1

# Update the "last updated" timestamp string
$ws.Range("A1").Value = "Datos actualizados a 15 de Abril de 2020 a las 21:22"
